$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Carga de archivos (si fuese necesario)" -> "Carga de archivos"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Carga de archivos (si fuese necesario)", $false, $false, $false, $false, $false, $true, 1, $false, "Carga de archivos", 2)

# ------------------------------------------------------------------
# Helper: find the paragraph index (1-based) whose trimmed text
# equals the given string. Returns -1 if not found.
# ------------------------------------------------------------------
function Get-ParaIndex($doc, $text) {
    $idx = -1
    $n = 1
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $text) {
            $idx = $n
        }
        $n = $n + 1
    }
    return $idx
}

# ------------------------------------------------------------------
# 2) Insert a new bulleted "firebase" item right after "jwt-decode"
#    (same list / numbering as that item).
# ------------------------------------------------------------------
$d = $word.ActiveDocument
$jwtIdx = Get-ParaIndex $d "jwt-decode"
$pJwt = $d.Paragraphs.Item($jwtIdx)
$pJwt.Range.InsertParagraphAfter()

$d = $word.ActiveDocument
$pFirebase = $d.Paragraphs.Item($jwtIdx + 1)
$pFirebase.Range.InsertAfter("firebase")

# ------------------------------------------------------------------
# 3) Insert three new bulleted items ("passport",
#    "passport-google-oauth20", "passport-facebook") right after
#    "jsonwebtoken" (same list / numbering as that item).
# ------------------------------------------------------------------
$d = $word.ActiveDocument
$jsonIdx = Get-ParaIndex $d "jsonwebtoken"

$pJson = $d.Paragraphs.Item($jsonIdx)
$pJson.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$pPassport = $d.Paragraphs.Item($jsonIdx + 1)
$pPassport.Range.InsertAfter("passport")

$d = $word.ActiveDocument
$pPassport = $d.Paragraphs.Item($jsonIdx + 1)
$pPassport.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$pGoogle = $d.Paragraphs.Item($jsonIdx + 2)
$pGoogle.Range.InsertAfter("passport-google-oauth20")

$d = $word.ActiveDocument
$pGoogle = $d.Paragraphs.Item($jsonIdx + 2)
$pGoogle.Range.InsertParagraphAfter()
$d = $word.ActiveDocument
$pFacebook = $d.Paragraphs.Item($jsonIdx + 3)
$pFacebook.Range.InsertAfter("passport-facebook")

# ------------------------------------------------------------------
# 4) Remove the last two blank paragraphs that sit directly above
#    the "Dependencias en el Backend" heading.
# ------------------------------------------------------------------
$d = $word.ActiveDocument
$backendIdx = Get-ParaIndex $d "Dependencias en el Backend"

$d.Paragraphs.Item($backendIdx - 1).Range.Delete()
$d = $word.ActiveDocument
$backendIdx = Get-ParaIndex $d "Dependencias en el Backend"
$d.Paragraphs.Item($backendIdx - 1).Range.Delete()
